$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1575.1578
$ws.Range("I40").Value = 1545.6923
$ws.Range("J40").Value = 1639
$ws.Range("K40").Value = 1545.6923
$ws.Range("L40").Value = 1639
$ws.Range("M40").Value = -1370.6923
$ws.Range("N40").Value = -1989

$ws.Range("H64").Value = 4117.7256
$ws.Range("I64").Value = 4068.9656
$ws.Range("J64").Value = 4160.5757
$ws.Range("K64").Value = 4068.9656
$ws.Range("L64").Value = 4160.5757
$ws.Range("M64").Value = -3820.9656
$ws.Range("N64").Value = -4656.5757

$ws.Range("H67").Value = 4117.7256
$ws.Range("I67").Value = 4068.9656
$ws.Range("J67").Value = 4160.5757
$ws.Range("K67").Value = 4068.9656
$ws.Range("L67").Value = 4160.5757
$ws.Range("M67").Value = -3210.9656
$ws.Range("N67").Value = -5876.5757

$ws.Range("H74").Value = 2582.8572
$ws.Range("I74").Value = 2245
$ws.Range("J74").Value = 3033.3333
$ws.Range("K74").Value = 2245
$ws.Range("L74").Value = 3033.3333
$ws.Range("M74").Value = -1309
$ws.Range("N74").Value = -4905.3333

$ws.Range("H77").Value = 2582.8572
$ws.Range("I77").Value = 2245
$ws.Range("J77").Value = 3033.3333
$ws.Range("K77").Value = 11225
$ws.Range("L77").Value = 15166.6665
$ws.Range("M77").Value = -6545
$ws.Range("N77").Value = -24526.6665

$ws.Range("H111").Value = 1608.5
$ws.Range("I111").Value = 1606.4445
$ws.Range("J111").Value = 1614.6666
$ws.Range("K111").Value = 4819.333500000001
$ws.Range("L111").Value = 4843.9998
$ws.Range("M111").Value = -1752.333500000001
$ws.Range("N111").Value = -10977.9998

$ws.Range("H113").Value = 1965.8948
$ws.Range("I113").Value = 1770.8
$ws.Range("J113").Value = 2035.5714
$ws.Range("K113").Value = 1770.8
$ws.Range("L113").Value = 2035.5714
$ws.Range("M113").Value = 1483.2
$ws.Range("N113").Value = -8543.571400000001

$ws.Range("H132").Value = 1877.3636
$ws.Range("I132").Value = 1765.1
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 5295.299999999999
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -2765.299999999999
$ws.Range("N132").Value = -14060

$ws.Range("H137").Value = 37471.355
$ws.Range("I137").Value = 1534.2632
$ws.Range("J137").Value = 113338.555
$ws.Range("K137").Value = 4602.7896
$ws.Range("L137").Value = 340015.665
$ws.Range("M137").Value = -2052.7896
$ws.Range("N137").Value = -345115.665

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2240.3333
$ws.Range("I2").Value = 2268.75
$ws.Range("J2").Value = 2013
$ws.Range("K2").Value = 2268.75
$ws.Range("L2").Value = 2013
$ws.Range("M2").Value = -2155.75
$ws.Range("N2").Value = -2239

$ws.Range("H88").Value = 2095.9167
$ws.Range("I88").Value = 2113.7273
$ws.Range("J88").Value = 1900
$ws.Range("K88").Value = 2113.7273
$ws.Range("L88").Value = 1900
$ws.Range("M88").Value = -1707.7273
$ws.Range("N88").Value = -2712

$ws.Range("H91").Value = 2095.9167
$ws.Range("I91").Value = 2113.7273
$ws.Range("J91").Value = 1900
$ws.Range("K91").Value = 2113.7273
$ws.Range("L91").Value = 1900
$ws.Range("M91").Value = -709.7273
$ws.Range("N91").Value = -4708

$ws.Range("H116").Value = 2240.3333
$ws.Range("I116").Value = 2268.75
$ws.Range("J116").Value = 2013
$ws.Range("K116").Value = 2268.75
$ws.Range("L116").Value = 2013
$ws.Range("M116").Value = 25.25
$ws.Range("N116").Value = -6601

$ws.Range("H132").Value = 3381.926
$ws.Range("I132").Value = 3156.25
$ws.Range("K132").Value = 9468.75
$ws.Range("M132").Value = -6938.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2240.3333
$ws.Range("I3").Value = 2268.75
$ws.Range("J3").Value = 2013
$ws.Range("K3").Value = 2268.75
$ws.Range("L3").Value = 2013
$ws.Range("M3").Value = -2154.75
$ws.Range("N3").Value = -2241

$ws.Range("H86").Value = 2629.625
$ws.Range("I86").Value = 2430.0588
$ws.Range("J86").Value = 3114.2856
$ws.Range("K86").Value = 2430.0588
$ws.Range("L86").Value = 3114.2856
$ws.Range("M86").Value = -1307.0588
$ws.Range("N86").Value = -5360.2856

$ws.Range("H89").Value = 2629.625
$ws.Range("I89").Value = 2430.0588
$ws.Range("J89").Value = 3114.2856
$ws.Range("K89").Value = 12150.294
$ws.Range("L89").Value = 15571.428
$ws.Range("M89").Value = -6534.293999999998
$ws.Range("N89").Value = -26803.428

$ws.Range("H134").Value = 2620.2273
$ws.Range("I134").Value = 2628.7334
$ws.Range("K134").Value = 7886.2002
$ws.Range("M134").Value = -5351.2002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 331.94116
$ws.Range("I22").Value = 339.7857
$ws.Range("J22").Value = 295.33334
$ws.Range("K22").Value = 339.7857
$ws.Range("L22").Value = 295.33334
$ws.Range("M22").Value = 10.21429999999998
$ws.Range("N22").Value = -995.33334

$ws.Range("H58").Value = 1864.3948
$ws.Range("I58").Value = 1331.1578
$ws.Range("J58").Value = 2397.6316
$ws.Range("K58").Value = 1331.1578
$ws.Range("L58").Value = 2397.6316
$ws.Range("M58").Value = -1128.1578
$ws.Range("N58").Value = -2803.6316

$ws.Range("H62").Value = 835957.3
$ws.Range("I62").Value = 1113643.1
$ws.Range("J62").Value = 2900
$ws.Range("K62").Value = 1113643.1
$ws.Range("L62").Value = 2900
$ws.Range("M62").Value = -1113019.1
$ws.Range("N62").Value = -4148

$ws.Range("H65").Value = 835957.3
$ws.Range("I65").Value = 1113643.1
$ws.Range("J65").Value = 2900
$ws.Range("K65").Value = 5568215.5
$ws.Range("L65").Value = 14500
$ws.Range("M65").Value = -5565095.5
$ws.Range("N65").Value = -20740

$ws.Range("H122").Value = 1264.6
$ws.Range("I122").Value = 1139.3158
$ws.Range("J122").Value = 1481
$ws.Range("K122").Value = 3417.9474
$ws.Range("L122").Value = 4443
$ws.Range("M122").Value = -967.9474
$ws.Range("N122").Value = -9343

$ws.Range("H134").Value = 3348.1667
$ws.Range("J134").Value = 7351.375
$ws.Range("L134").Value = 22054.125
$ws.Range("N134").Value = -27124.125

$ws.Range("H136").Value = 1864.3948
$ws.Range("I136").Value = 1331.1578
$ws.Range("J136").Value = 2397.6316
$ws.Range("K136").Value = 3993.4734
$ws.Range("L136").Value = 7192.8948
$ws.Range("M136").Value = -1443.4734
$ws.Range("N136").Value = -12292.8948

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 66725.60000000001
$ws.Range("I2").Value = 125027.5
$ws.Range("J2").Value = 94.85714
$ws.Range("K2").Value = 750165
$ws.Range("L2").Value = 569.14284
$ws.Range("M2").Value = -750052
$ws.Range("N2").Value = -795.14284

$ws.Range("H23").Value = 149.41667
$ws.Range("I23").Value = 90
$ws.Range("J23").Value = 161.3
$ws.Range("K23").Value = 270
$ws.Range("L23").Value = 483.9
$ws.Range("M23").Value = -35
$ws.Range("N23").Value = -953.9000000000001

$ws.Range("H103").Value = 520.6667
$ws.Range("I103").Value = 422
$ws.Range("J103").Value = 570
$ws.Range("K103").Value = 1266
$ws.Range("L103").Value = 1710
$ws.Range("M103").Value = -387
$ws.Range("N103").Value = -3468

$ws.Range("H131").Value = 738.39703
$ws.Range("I131").Value = 323.43478
$ws.Range("J131").Value = 950.4888999999999
$ws.Range("K131").Value = 970.3043399999999
$ws.Range("L131").Value = 2851.4667
$ws.Range("M131").Value = 4069.69566
$ws.Range("N131").Value = -12931.4667

$ws.Range("H132").Value = 230255.66
$ws.Range("I132").Value = 504.87097
$ws.Range("K132").Value = 4543.83873
$ws.Range("M132").Value = -2013.83873

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 100.56
$ws.Range("I2").Value = 55.454544
$ws.Range("J2").Value = 136
$ws.Range("K2").Value = 55.454544
$ws.Range("L2").Value = 136
$ws.Range("M2").Value = 57.545456
$ws.Range("N2").Value = -362

$ws.Range("H70").Value = 6600
$ws.Range("I70").Value = 7794.1177
$ws.Range("J70").Value = 4570
$ws.Range("K70").Value = 7794.1177
$ws.Range("L70").Value = 4570
$ws.Range("M70").Value = -7524.1177
$ws.Range("N70").Value = -5110

$ws.Range("H73").Value = 6600
$ws.Range("I73").Value = 7794.1177
$ws.Range("J73").Value = 4570
$ws.Range("K73").Value = 7794.1177
$ws.Range("L73").Value = 4570
$ws.Range("M73").Value = -6858.1177
$ws.Range("N73").Value = -6442

$ws.Range("H126").Value = 2842.2144
$ws.Range("I126").Value = 2316.353
$ws.Range("J126").Value = 3654.9092
$ws.Range("K126").Value = 6949.059
$ws.Range("L126").Value = 10964.7276
$ws.Range("M126").Value = -4479.059
$ws.Range("N126").Value = -15904.7276

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2050.0908
$ws.Range("I7").Value = 1625.3334
$ws.Range("J7").Value = 2559.8
$ws.Range("K7").Value = 1625.3334
$ws.Range("L7").Value = 2559.8
$ws.Range("M7").Value = -1513.3334
$ws.Range("N7").Value = -2783.8

$ws.Range("H61").Value = 12894
$ws.Range("I61").Value = 15872.286
$ws.Range("J61").Value = 2470
$ws.Range("K61").Value = 15872.286
$ws.Range("L61").Value = 2470
$ws.Range("M61").Value = -15670.286
$ws.Range("N61").Value = -2874

$ws.Range("H113").Value = 12894
$ws.Range("I113").Value = 15872.286
$ws.Range("J113").Value = 2470
$ws.Range("K113").Value = 15872.286
$ws.Range("L113").Value = 2470
$ws.Range("M113").Value = -13702.286
$ws.Range("N113").Value = -6810

$ws.Range("H126").Value = 2050.0908
$ws.Range("I126").Value = 1625.3334
$ws.Range("J126").Value = 2559.8
$ws.Range("K126").Value = 4876.0002
$ws.Range("L126").Value = 7679.400000000001
$ws.Range("M126").Value = -2406.0002
$ws.Range("N126").Value = -12619.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 706.7
$ws.Range("I107").Value = 599.2
$ws.Range("J107").Value = 814.2
$ws.Range("K107").Value = 1797.6
$ws.Range("L107").Value = 2442.6
$ws.Range("M107").Value = 122.3999999999999
$ws.Range("N107").Value = -6282.6

$ws.Range("H126").Value = 1156.909
$ws.Range("I126").Value = 989.5714
$ws.Range("J126").Value = 1449.75
$ws.Range("K126").Value = 2968.7142
$ws.Range("L126").Value = 4349.25
$ws.Range("M126").Value = -498.7142000000003
$ws.Range("N126").Value = -9289.25
